# Updated cryptos list on Tue Jun 25 22:01:16 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures pulled from coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number (e.g. '577.27') need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric values and
# the original text formatting (fixed decimals, thousand-dot grouping) would be lost.
$textPriceCells = @("D5", "D6", "D9", "D10", "D16", "D21", "D22", "D26", "D29", "D30", "D35", "D39", "D40", "D42", "D46", "D47", "D48", "D49")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.895.06'
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("D3").Value = '3.401.60'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '577.27'
$ws.Range("E5").Value = '  +2.33%  '
$ws.Range("D6").Value = '137.37'
$ws.Range("E6").Value = '  +5.02%  '
$ws.Range("E8").Value = '  +0.85%  '
$ws.Range("D9").Value = '7.51'
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("D10").Value = '0.126'
$ws.Range("E10").Value = '  +7.20%  '
$ws.Range("E11").Value = '  +4.30%  '
$ws.Range("D12").Value = '3.984.05'
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("E14").Value = '  +5.71%  '
$ws.Range("D15").Value = '3.415.14'
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").Value = '25.45'
$ws.Range("E16").Value = '  +3.04%  '
$ws.Range("D17").Value = '61.946.86'
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("E18").Value = '  +6.01%  '
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("E20").Value = '  +5.11%  '
$ws.Range("D21").Value = '389.37'
$ws.Range("E21").Value = '  +9.92%  '
$ws.Range("D22").Value = '0.571'
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("D23").Value = '3.543.92'
$ws.Range("E23").Value = '  +1.95%  '
$ws.Range("E24").Value = '  +15.21%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '71.51'
$ws.Range("E26").Value = '  +3.16%  '
$ws.Range("E27").Value = '  +2.87%  '
$ws.Range("E28").Value = '  -4.72%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '8.29'
$ws.Range("E30").Value = '  +4.49%  '
$ws.Range("E31").Value = '  +3.69%  '
$ws.Range("E32").Value = '  +2.41%  '
$ws.Range("D34").Value = '3.433.52'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("D35").Value = '23.57'
$ws.Range("E35").Value = '  +2.77%  '
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("E38").Value = '  +3.80%  '
$ws.Range("D39").Value = '164.09'
$ws.Range("E39").Value = '  +4.29%  '
$ws.Range("D40").Value = '0.0788'
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("E41").Value = '  +13.44%  '
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  +2.83%  '
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("D46").Value = '25.12'
$ws.Range("E46").Value = '  +5.91%  '
$ws.Range("D47").Value = '41.65'
$ws.Range("E47").Value = '  +2.37%  '
$ws.Range("D48").Value = '7.00'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").Value = '23.29'
$ws.Range("E49").Value = '  +3.78%  '
$ws.Range("D50").Value = '2.367.77'
$ws.Range("E50").Value = '  +8.73%  '
$ws.Range("E51").Value = '  +6.61%  '
